$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13-21 down to 14-22),
# matching the new "Docentes responsaveis" value row placed right after the
# "Objectives:"/"Docentes responsaveis:" label rows (11-12).
$ws.Rows.Item(13).Insert()

# Column A labels already line up correctly after the shift, only the B/C
# value columns require updating (old values were duplicated/misaligned).
$ws.Cells.Item(10, 2).Value = "Apresentar o conjunto de métodos de investigação de superfície e de sub-superfície (diretos e indiretos) utilizados nas caracterizações geológico-geotécnicas que envolvem o meio ambiente. Estabelecer análise crítica que possibilite a escolha e a utilização adequadas das técnicas de investigações disponíveis visando o estudo dos diversos tipos de problemas ambientais."
$ws.Cells.Item(10, 3).Value = "Apresentar o conjunto de métodos de investigação de superfície e de sub-superfície (diretos e indiretos) utilizados nas caracterizações geológico-geotécnicas que envolvem o meio ambiente. Estabelecer análise crítica que possibilite a escolha e a utilização adequadas das técnicas de investigações disponíveis visando o estudo dos diversos tipos de problemas ambientais."
$ws.Cells.Item(13, 2).Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Cells.Item(13, 3).Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Cells.Item(14, 2).Value = "Investigação de superfície e sub-superfícies; Técnicas e equipamentos mais adequados."
$ws.Cells.Item(14, 3).Value = "Investigação de superfície e sub-superfícies; Técnicas e equipamentos mais adequados."
$ws.Cells.Item(16, 2).Value = "Introdução, conceitos e objetivos; Seqüência de estudos rotineiros (usuais); Estudos Corretivos e Preventivos; Investigação de Superfície; Investigação de Sub-superfície - Aplicações / Limitações - Métodos Diretos; Métodos Indiretos (Geofísicos); Métodos Sísmicos; Métodos Elétricos e Eletromagnético; Ensaios em Furos de Sondagem; Ensaios com Traçadores; Instrumentação Hidráulica e Mecânica. Estudo de caso."
$ws.Cells.Item(16, 3).Value = "Introdução, conceitos e objetivos; Seqüência de estudos rotineiros (usuais); Estudos Corretivos e Preventivos; Investigação de Superfície; Investigação de Sub-superfície - Aplicações / Limitações - Métodos Diretos; Métodos Indiretos (Geofísicos); Métodos Sísmicos; Métodos Elétricos e Eletromagnético; Ensaios em Furos de Sondagem; Ensaios com Traçadores; Instrumentação Hidráulica e Mecânica. Estudo de caso."
$ws.Cells.Item(19, 2).Value = "Aulas expositivas, exercícios e visitas didátias de campo."
$ws.Cells.Item(19, 3).Value = "Aulas expositivas, exercícios e visitas didátias de campo."
$ws.Cells.Item(20, 2).Value = "Provas e relatórios."
$ws.Cells.Item(20, 3).Value = "Provas e relatórios."
$ws.Cells.Item(21, 2).Value = "Prova única com nota igual ou superior a 5,0 (cinco)."
$ws.Cells.Item(21, 3).Value = "Prova única com nota igual ou superior a 5,0 (cinco)."
$ws.Cells.Item(22, 2).Value = "ATTEWELL & FARMER - 1976 - Principles of Enginnering Geology. Chapman Hall.`nDUNICLIFF, J. - 1988 - Geotechnical Instrumentation for Monitoring Field Performance, Joh Willey & Sons, New York, 577 p.`nHANNA, T.H. - 1996 - Field Instrumentation in Geotechnical Engineering. Trans Tech Publications, RockPort - MA, 843 p.`nKELLY, W.E. e MARES S. - Applied Geophyses in Hydrogeological and Engineering Practice. Elsevier, New York - 1993, 300p.`nLUIZ, J.G. - 1995 - Geofísica de Prospecção. Editora Universitária UFPA, Belém, 1995. `nVOGELSAND, D. - 1995 - Environmental Geophysics. Springer - Verlag, Berlin, 171p."
$ws.Cells.Item(22, 3).Value = "ATTEWELL & FARMER - 1976 - Principles of Enginnering Geology. Chapman Hall.`nDUNICLIFF, J. - 1988 - Geotechnical Instrumentation for Monitoring Field Performance, Joh Willey & Sons, New York, 577 p.`nHANNA, T.H. - 1996 - Field Instrumentation in Geotechnical Engineering. Trans Tech Publications, RockPort - MA, 843 p.`nKELLY, W.E. e MARES S. - Applied Geophyses in Hydrogeological and Engineering Practice. Elsevier, New York - 1993, 300p.`nLUIZ, J.G. - 1995 - Geofísica de Prospecção. Editora Universitária UFPA, Belém, 1995. `nVOGELSAND, D. - 1995 - Environmental Geophysics. Springer - Verlag, Berlin, 171p."
